$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.6184

$ws.Range("A4").Value = -21.36130000000002
$ws.Range("B4").Value = 5.361199999999994
$ws.Range("E4").Value = 12.99170000000002

$ws.Range("B5").Value = 4.926800000000003

$ws.Range("A6").Value = -21.13010000000001
$ws.Range("B6").Value = 5.992799999999993

$ws.Range("A7").Value = -21.08620000000001

$ws.Range("A8").Value = -21.39120000000002
$ws.Range("B8").Value = 5.359799999999995

$ws.Range("E9").Value = 13.91940000000001

$ws.Range("E11").Value = 13.79680000000001

$ws.Range("E14").Value = 12.5756

$ws.Range("A16").Value = -21.30020000000002
$ws.Range("B16").Value = 5.368599999999998

$ws.Range("E18").Value = 13.17980000000001

$ws.Range("A20").Value = -22.07240000000003

$ws.Range("A21").Value = -20.20429999999999

$ws.Range("B22").Value = 4.938700000000003

$ws.Range("E25").Value = 12.58000000000001
